# "hero levelup & hero refine"
#
# The separate LevelupExp sheet is removed; HeroProto instead grows a new
# "晋级材料" (upgrade materials) column (K) that points at the per-level
# material-needs maps.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("LevelupExp").Delete()

$ws = $wb.Worksheets.Item("HeroProto")
$ws.Activate()

# --- Column K header block (rows 1-4), mirroring the other columns'
#     field-name / type / comment / map-type header rows. ---
$ws.Range("K1").Value = "晋级材料"
$ws.Range("K2").Value = "jl"
# Give this "jl" its own distinct entry in the shared-string table (like
# several of the sheet's other short field codes already have) rather
# than silently collapsing onto an unrelated existing "jl" cell.
$ws.Range("K2").Characters(1, 1).Font.Name = "宋体"
$ws.Range("K3").Value = "needs#id_cnt"
$ws.Range("K4").Value = "map"
# Likewise keep "map" as its own entry.
$ws.Range("K4").Characters(1, 1).Font.Name = "宋体"

# --- Data rows: each hero/level row references the 4113_x material map. ---
$materials = "4113_1", "4113_2", "4113_3", "4113_4", "4113_5"
for ($row = 5; $row -le 29; $row++) {
    $ws.Range("K$row").Value = $materials[($row - 5) % 5]
}

# Match the existing header/data cell styling used throughout the sheet
# (font 4, same as columns A/B/E/H/I/J) by copying format from a cell that
# already carries it.
$ws.Range("A2").Copy()
$ws.Range("K1:K29").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Size the new column like the sheet's other bestFit text columns.
$ws.Columns.Item(11).ColumnWidth = 12

# Leave the selection where the author ended up.
$ws.Range("K13").Select()
